$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
  "2+49=51",
  "4+62=66",
  "4+23=27",
  "61-46=15",
  "14+44=58",
  "96-44=52",
  "0+9=9",
  "66-10=56",
  "72+7=79",
  "24+31=55",
  "48-1=47",
  "81-15=66",
  "84-46=38",
  "86+13=99",
  "79-12=67",
  "58-36=22",
  "14+0=14",
  "85-59=26",
  "53+34=87",
  "56+0=56",
  "24+0=24",
  "73+22=95",
  "37-22=15",
  "84-61=23",
  "61+16=77",
  "79+10=89",
  "12+50=62",
  "29+59=88",
  "86-26=60",
  "98-82=16",
  "28+1=29",
  "5+11=16",
  "66-18=48",
  "24+55=79",
  "50-48=2",
  "2+7=9",
  "97-74=23",
  "88-67=21",
  "95-31=64",
  "69-5=64",
  "43-2=41",
  "82-8=74",
  "6+92=98",
  "27+29=56",
  "24-3=21",
  "45+14=59",
  "44+48=92",
  "13+6=19",
  "76-18=58",
  "15+30=45",
  "50-49=1",
  "70+2=72",
  "85-43=42",
  "14+44=58",
  "12+45=57",
  "8+78=86",
  "82-21=61",
  "3+16=19",
  "93-85=8",
  "11+19=30",
  "10-7=3",
  "3+13=16",
  "89-36=53",
  "19+36=55",
  "13-2=11",
  "37+2=39",
  "98-9=89",
  "47-36=11",
  "51+17=68",
  "83-57=26",
  "94-84=10",
  "55-39=16",
  "17-2=15",
  "87-26=61",
  "0+73=73",
  "68+19=87",
  "72-14=58",
  "66-26=40",
  "87-68=19",
  "77-42=35",
  "1+26=27",
  "64+30=94",
  "24+28=52",
  "66-27=39",
  "86-49=37",
  "63+27=90",
  "65-52=13",
  "33+22=55",
  "56-34=22",
  "10+0=10",
  "8+6=14",
  "66-36=30",
  "53-49=4",
  "0+95=95",
  "53-41=12",
  "58-24=34",
  "46-12=34",
  "32-8=24",
  "32+21=53",
  "2+60=62"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    if ($idx -lt $newValues.Length) {
      $cell = $t.Cell($r, $c)
      $cell.Range.Text = $newValues[$idx]
    }
    $idx = $idx + 1
  }
}

Write-Host "Updated" $idx "cells"
